$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for the season record: Wins, Losses, Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, centered, bordered) by copying
# the format from the adjacent header cell instead of re-creating a style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in each player's row with the team's season record (97-65-0)
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 97  # column AD - Wins
    $ws.Cells.Item($r, 31).Value = 65  # column AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # column AF - Ties
}
